# Auto-generated: re-apply market-price snapshot values from the scheduled runner.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 44070.332
$ws.Range("J95").Value = 44070.332
$ws.Range("L95").Value = 44070.332
$ws.Range("N95").Value = -49562.332
$ws.Range("H116").Value = 7001.3335
$ws.Range("I116").Value = 7001.3335
$ws.Range("K116").Value = 7001.3335
$ws.Range("M116").Value = -3559.3335
$ws.Range("H129").Value = 3031907.5
$ws.Range("I129").Value = 33333764
$ws.Range("K129").Value = 100001292
$ws.Range("M129").Value = -99996292

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 644
$ws.Range("I2").Value = 642
$ws.Range("J2").Value = 650
$ws.Range("K2").Value = 642
$ws.Range("L2").Value = 650
$ws.Range("M2").Value = -529
$ws.Range("N2").Value = -876
$ws.Range("H63").Value = 6400.625
$ws.Range("I63").Value = 6400.625
$ws.Range("K63").Value = 6400.625
$ws.Range("M63").Value = -5714.625
$ws.Range("H66").Value = 6400.625
$ws.Range("I66").Value = 6400.625
$ws.Range("K66").Value = 32003.125
$ws.Range("M66").Value = -28571.125
$ws.Range("H97").Value = 189.8
$ws.Range("I97").Value = 189.8
$ws.Range("K97").Value = 189.8
$ws.Range("M97").Value = 306.2
$ws.Range("H110").Value = 697
$ws.Range("I110").Value = 697
$ws.Range("K110").Value = 697
$ws.Range("M110").Value = 1348
$ws.Range("H116").Value = 644
$ws.Range("I116").Value = 642
$ws.Range("J116").Value = 650
$ws.Range("K116").Value = 642
$ws.Range("L116").Value = 650
$ws.Range("M116").Value = 1652
$ws.Range("N116").Value = -5238
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 644
$ws.Range("I3").Value = 642
$ws.Range("J3").Value = 650
$ws.Range("K3").Value = 642
$ws.Range("L3").Value = 650
$ws.Range("M3").Value = -528
$ws.Range("N3").Value = -878

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 3461.125
$ws.Range("I11").Value = 922.25
$ws.Range("K11").Value = 2766.75
$ws.Range("M11").Value = -2626.75
$ws.Range("H114").Value = 315.5
$ws.Range("J114").Value = 359.66666
$ws.Range("L114").Value = 1078.99998
$ws.Range("N114").Value = -7586.999980000001
$ws.Range("H140").Value = 5498.25
$ws.Range("I140").Value = 5498.25
$ws.Range("K140").Value = 16494.75
$ws.Range("M140").Value = -11314.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 800
$ws.Range("I70").Value = 800
$ws.Range("K70").Value = 800
$ws.Range("M70").Value = -530
$ws.Range("H73").Value = 800
$ws.Range("I73").Value = 800
$ws.Range("K73").Value = 800
$ws.Range("M73").Value = 136
$ws.Range("H80").Value = 2370
$ws.Range("I80").Value = 2137.5
$ws.Range("K80").Value = 2137.5
$ws.Range("M80").Value = -1139.5
$ws.Range("H83").Value = 2370
$ws.Range("I83").Value = 2137.5
$ws.Range("K83").Value = 10687.5
$ws.Range("M83").Value = -5695.5
$ws.Range("H113").Value = 2200.9
$ws.Range("I113").Value = 2223.2222
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 2223.2222
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -53.22220000000016
$ws.Range("N113").Value = -6340
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2051.7334
$ws.Range("I16").Value = 2250.1667
$ws.Range("J16").Value = 1919.4445
$ws.Range("K16").Value = 2250.1667
$ws.Range("L16").Value = 1919.4445
$ws.Range("M16").Value = -2080.1667
$ws.Range("N16").Value = -2259.4445
$ws.Range("H60").Value = 34997.5
$ws.Range("J60").Value = 34997.5
$ws.Range("L60").Value = 34997.5
$ws.Range("N60").Value = -36015.5
$ws.Range("H68").Value = 2740
$ws.Range("I68").Value = 3100
$ws.Range("K68").Value = 3100
$ws.Range("M68").Value = -2351
$ws.Range("H71").Value = 2740
$ws.Range("I71").Value = 3100
$ws.Range("K71").Value = 15500
$ws.Range("M71").Value = -11756

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 10335
$ws.Range("J5").Value = 10335
$ws.Range("L5").Value = 10335
$ws.Range("N5").Value = -10559
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H68").Value = 37492.5
$ws.Range("J68").Value = 37492.5
$ws.Range("L68").Value = 37492.5
$ws.Range("N68").Value = -39114.5
$ws.Range("H71").Value = 37492.5
$ws.Range("J71").Value = 37492.5
$ws.Range("L71").Value = 112477.5
$ws.Range("N71").Value = -120589.5
$ws.Range("H81").Value = 750
$ws.Range("I81").Value = 1000
$ws.Range("J81").Value = 500
$ws.Range("K81").Value = 2000
$ws.Range("L81").Value = 1000
$ws.Range("M81").Value = -939
$ws.Range("N81").Value = -3122
$ws.Range("H84").Value = 750
$ws.Range("I84").Value = 1000
$ws.Range("J84").Value = 500
$ws.Range("K84").Value = 10000
$ws.Range("L84").Value = 5000
$ws.Range("M84").Value = -4696
$ws.Range("N84").Value = -15608
$ws.Range("H97").Value = 20572
$ws.Range("J97").Value = 20572
$ws.Range("L97").Value = 20572
$ws.Range("N97").Value = -22554
$ws.Range("H113").Value = 1873.091
$ws.Range("I113").Value = 199.5
$ws.Range("J113").Value = 2245
$ws.Range("K113").Value = 598.5
$ws.Range("L113").Value = 6735
$ws.Range("M113").Value = 1571.5
$ws.Range("N113").Value = -11075
